{"js": "// Apply the \"Setting Up NGons Tutorial\" edits to the SEO meta document.\n// 1. Title: \"1 What is Blender\" -> \"The Knife Tool\"\n// 2. Keywords: insert \"The Knife Tool, \" before \"Blender, 3D Modeling, Animation, Graphic Art\"\n// 3. Description: replace the \"what the 3D modeling program ... about.\" sentence with the\n//    new Knife Tool description (keeping the trailing \"/>\" literal).\n// 4. Category: insert \"The Knife Tool, \" before \"Blender, 3D Modeling, Animation, Graphic Art\"\n// 5. Revised date: \"Wednesday, December 11, 2024\" -> \"Sunday, January 12, 2025\"\n// 6. Url: point at the new 2025 Knife Tool article location.\n\nconst body = context.document.body;\n\n// --- 1. Title -------------------------------------------------------------\nconst titleResults = body.search(\"1 What is Blender\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"The Knife Tool\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2 & 4. Keywords + Category: insert \"The Knife Tool, \" before each\n//     occurrence of \"Blender, 3D Modeling, Animation, Graphic Art\"\nconst kwCatResults = body.search(\"Blender, 3D Modeling, Animation, Graphic Art\", { matchCase: true });\nkwCatResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < kwCatResults.items.length; i++) {\n  kwCatResults.items[i].insertText(\"The Knife Tool, \", Word.InsertLocation.before);\n}\nawait context.sync();\n\n// --- 3. Description --------------------------------------------------------\n// Locate the paragraph that contains the word \"description\" (the <meta name=\"description\" .../>\n// line) and operate within it so we don't disturb any other paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet descParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"what the 3D modeling program\") !== -1) {\n    descParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (descParagraph) {\n  const startIdx = descParagraph.text.indexOf(\"what the 3D modeling program\");\n  const tail = descParagraph.text.substring(startIdx); // everything through end of paragraph (\"...about./>\")\n\n  const descResults = descParagraph.search(tail, { matchCase: true });\n  descResults.load(\"items\");\n  await context.sync();\n\n  if (descResults.items.length > 0) {\n    // Remove the whole old sentence (this also clears out the now-orphaned\n    // proofErr grammar markers around \"about.\").\n    descResults.items[0].delete();\n    await context.sync();\n\n    // Re-insert the replacement text as two runs, matching the target shape:\n    // \"how to use the Knife tool, while in Edit mode inside of the Blender application\" + \"/>\"\n    descParagraph.insertText(\n      \"how to use the Knife tool, while in Edit mode inside of the Blender application\",\n      Word.InsertLocation.end\n    );\n    await context.sync();\n    descParagraph.insertText(\"/>\", Word.InsertLocation.end);\n    await context.sync();\n  }\n}\n\n// --- 5. Revised date --------------------------------------------------------\nconst dateResults = body.search(\"Wednesday, December 11, 2024\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"Sunday, January 12, 2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 6. Url -----------------------------------------------------------------\nconst urlResults = body.search(\n  \"Enlightenment/Articles/2024/8-Blender-2024/1-What-Is-Blender/1-What-Is-Blender.html\",\n  { matchCase: true }\n);\nurlResults.load(\"items\");\nawait context.sync();\nif (urlResults.items.length > 0) {\n  urlResults.items[0].insertText(\n    \"Enlightenment/Articles/2025/1-Blender-Continued/2-Edit-Mode/1-The-Menus/1-The-Tools-Menu/5-Knife-Tool/The-Knife-Tool.html\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Apply the \"Setting Up NGons Tutorial\" edits to the SEO meta document.\n# 1. Title: \"1 What is Blender\" -> \"The Knife Tool\"\n# 2. Keywords: insert \"The Knife Tool, \" before \"Blender, 3D Modeling, Animation, Graphic Art\"\n# 3. Description: replace the \"what the 3D modeling program ... about.\" sentence with the\n#    new Knife Tool description (keeping the trailing \"/>\" literal).\n# 4. Category: insert \"The Knife Tool, \" before \"Blender, 3D Modeling, Animation, Graphic Art\"\n# 5. Revised date: \"Wednesday, December 11, 2024\" -> \"Sunday, January 12, 2025\"\n# 6. Url: point at the new 2025 Knife Tool article location.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title ---------------------------------------------------------------\n$titleRange = $d.Content\nif ($titleRange.Find.Execute(\"1 What is Blender\")) {\n  $titleRange.Text = \"The Knife Tool\"\n}\n\n# --- 2 & 4. Keywords + Category: insert \"The Knife Tool, \" before each\n#     occurrence of \"Blender, 3D Modeling, Animation, Graphic Art\"\n$kwCatRange = $d.Content\nwhile ($kwCatRange.Find.Execute(\"Blender, 3D Modeling, Animation, Graphic Art\")) {\n  $insertPoint = $d.Range($kwCatRange.Start, $kwCatRange.Start)\n  $insertPoint.InsertBefore(\"The Knife Tool, \")\n  $kwCatRange.Collapse(0)\n}\n\n# --- 3. Description -----------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  $pText = $p.Range.Text\n  if ($pText -like \"*what the 3D modeling program*\") {\n    $pRange = $p.Range\n    $found = $pRange.Find.Execute(\"what the 3D modeling program\")\n    $startPos = $pRange.Start\n    $paraEnd = $p.Range.End - 1  # exclude the paragraph mark\n\n    # Remove the old sentence (also clears out the now-orphaned proofErr\n    # grammar markers that surrounded \"about.\").\n    $targetRange = $d.Range($startPos, $paraEnd)\n    $targetRange.Delete()\n\n    # Re-insert the replacement text as two runs, matching the target shape:\n    # \"how to use the Knife tool, while in Edit mode inside of the Blender application\" + \"/>\"\n    $insertRange = $d.Range($startPos, $startPos)\n    $insertRange.InsertAfter(\"how to use the Knife tool, while in Edit mode inside of the Blender application\")\n    $insertRange.Collapse(0)\n    $insertRange.InsertAfter(\"/>\")\n  }\n}\n\n# --- 5. Revised date ----------------------------------------------------------\n$dateRange = $d.Content\nif ($dateRange.Find.Execute(\"Wednesday, December 11, 2024\")) {\n  $dateRange.Text = \"Sunday, January 12, 2025\"\n}\n\n# --- 6. Url ---------------------------------------------------------------------\n$urlRange = $d.Content\nif ($urlRange.Find.Execute(\"Enlightenment/Articles/2024/8-Blender-2024/1-What-Is-Blender/1-What-Is-Blender.html\")) {\n  $urlRange.Text = \"Enlightenment/Articles/2025/1-Blender-Continued/2-Edit-Mode/1-The-Menus/1-The-Tools-Menu/5-Knife-Tool/The-Knife-Tool.html\"\n}\n"}
